# Weekly update: a new price record for the week is inserted at the top of
# the "Feria Lagunitas de Puerto Montt - Apio" data block (row 186), pushing
# all the existing records (previously rows 186-219) down by one row
# (now rows 187-220), and the sheet's used range grows from A1:R219 to
# A1:R220.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 186; this shifts rows 186..219 down to 187..220,
# carrying their values/formatting with them (matches Excel's native
# "insert row" behaviour).
$ws.Rows.Item(186).Insert()

# Populate the newly-inserted row 186 with this week's record.
$ws.Cells.Item(186, 1).Value  = 4
$ws.Cells.Item(186, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(186, 3).Value  = "Los Lagos"
$ws.Cells.Item(186, 4).Value  = 44637
$ws.Cells.Item(186, 5).Value  = 10
$ws.Cells.Item(186, 6).Value  = 100112017
$ws.Cells.Item(186, 7).Value  = "Apio"
$ws.Cells.Item(186, 8).Value  = "Americana (o)"
$ws.Cells.Item(186, 9).Value  = "Primera"
$ws.Cells.Item(186, 10).Value = 35
$ws.Cells.Item(186, 11).Value = 13000
$ws.Cells.Item(186, 12).Value = 13000
$ws.Cells.Item(186, 13).Value = 13000
$ws.Cells.Item(186, 14).Value = "$/docena de matas"
$ws.Cells.Item(186, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(186, 16).Value = 2167
$ws.Cells.Item(186, 17).Value = 6
$ws.Cells.Item(186, 18).Value = "Hortaliza"
